$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column C - copy the header formatting from B1 (same style as A1)
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C1").Value = "ordem"

# "ordem" values per row - sequential order number per distinct id group in column A
$ordem = @{
    2 = 1; 3 = 1; 4 = 1; 5 = 1; 6 = 1; 7 = 1; 8 = 1; 9 = 1; 10 = 1
    11 = 2; 12 = 2; 13 = 2; 14 = 2
    15 = 3; 16 = 3; 17 = 3
    18 = 4; 19 = 4
    20 = 5; 21 = 5; 22 = 5
    23 = 6; 24 = 6; 25 = 6; 26 = 6; 27 = 6; 28 = 6
    29 = 7; 30 = 7
    31 = 8; 32 = 8
    33 = 9; 34 = 9; 35 = 9; 36 = 9; 37 = 9
    38 = 10; 39 = 10; 40 = 10; 41 = 10
}

foreach ($row in $ordem.Keys) {
    $ws.Cells.Item($row, 3).Value = $ordem[$row]
}
